$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chargingdata")

# Insert three new rows at position 5 (shifting existing rows 5.. down to 8..)
$ws.Range("A5:A7").EntireRow.Insert()

# Populate the newly inserted rows 5-7 with copies of the data from rows 2-4
# (the 2018-06 / 2018-07 / 2018-08 rows), matching the layout added upstream.
$ws.Range("A2:F2").Copy($ws.Range("A5:F5"))
$ws.Range("A3:F3").Copy($ws.Range("A6:F6"))
$ws.Range("A4:F4").Copy($ws.Range("A7:F7"))

# Restore the selected cell to I8, as recorded in the saved workbook view.
$ws.Range("I8").Select()
